# Weekly crime data update (cs-en-us-025pct)
# - Bulletin header: Volume 32, Number 11 -> 12
# - Report week: 3/10/2025-3/16/2025 -> 3/17/2025-3/23/2025
# - Updated weekly crime complaint counts/percentages for rows 14-30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/17/2025  Through  3/23/2025"

# --- Crime Complaints table updates (rows 14-30) ---
# Cells that previously held the "N/A" / "***.*" placeholder text now hold
# real numbers, so their number format is copied from a same-column sibling
# cell that already uses the numeric style.
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = $ws.Range("C15").NumberFormat
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = $ws.Range("K14").NumberFormat
$ws.Range("G14").Value = 1
$ws.Range("G14").NumberFormat = $ws.Range("C15").NumberFormat
$ws.Range("H14").Value = -100
$ws.Range("H14").NumberFormat = $ws.Range("K14").NumberFormat
$ws.Range("J14").Value = 2
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 200
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = 100
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -26.315789473684
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = -25.423728813559
$ws.Range("L16").Value = -16.981132075471
$ws.Range("M16").Value = -8.333333333333
$ws.Range("N16").Value = -73.006134969325
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -61.538461538461
$ws.Range("I17").Value = 56
$ws.Range("J17").Value = 108
$ws.Range("K17").Value = -48.148148148148
$ws.Range("L17").Value = -30
$ws.Range("M17").Value = 19.148936170212
$ws.Range("N17").Value = -52.941176470588
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("I18").Value = 25
$ws.Range("J18").Value = 19
$ws.Range("K18").Value = 31.578947368421
$ws.Range("M18").Value = 13.636363636363
$ws.Range("N18").Value = -75.961538461538
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -27.659574468085
$ws.Range("I19").Value = 83
$ws.Range("J19").Value = 101
$ws.Range("K19").Value = -17.821782178217
$ws.Range("L19").Value = -4.597701149425
$ws.Range("M19").Value = 97.619047619047
$ws.Range("N19").Value = 56.603773584905
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -80
$ws.Range("I20").Value = 12
$ws.Range("J20").Value = 24
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = 9.090909090909
$ws.Range("M20").Value = 33.333333333333
$ws.Range("N20").Value = -82.089552238806
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -31.578947368421
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 133
$ws.Range("H21").Value = -43.609022556391
$ws.Range("I21").Value = 226
$ws.Range("J21").Value = 316
$ws.Range("K21").Value = -28.481012658227
$ws.Range("L21").Value = -12.403100775193
$ws.Range("M21").Value = 29.142857142857
$ws.Range("N21").Value = -56.031128404669
$ws.Range("D22").Value = 2
$ws.Range("D22").NumberFormat = $ws.Range("C15").NumberFormat
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = $ws.Range("K14").NumberFormat
$ws.Range("G22").Value = 3
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -71.428571428571
$ws.Range("L22").Value = -60
$ws.Range("M22").Value = -83.333333333333
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 34
$ws.Range("J23").Value = 43
$ws.Range("K23").Value = -20.930232558139
$ws.Range("L23").Value = -8.108108108108
$ws.Range("M23").Value = 100
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -67.857142857142
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = -30.769230769230
$ws.Range("I24").Value = 172
$ws.Range("J24").Value = 201
$ws.Range("K24").Value = -14.427860696517
$ws.Range("L24").Value = -42.857142857142
$ws.Range("M24").Value = -20.737327188940
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -17.647058823529
$ws.Range("I25").Value = 47
$ws.Range("J25").Value = 49
$ws.Range("K25").Value = -4.081632653061
$ws.Range("L25").Value = -67.132867132867
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -6.666666666666
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 72
$ws.Range("H26").Value = -44.444444444444
$ws.Range("I26").Value = 131
$ws.Range("J26").Value = 149
$ws.Range("K26").Value = -12.080536912751
$ws.Range("L26").Value = 18.018018018018
$ws.Range("M26").Value = 36.458333333333
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = 0
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 300
$ws.Range("I28").Value = 10
$ws.Range("K28").Value = -9.090909090909
$ws.Range("L28").Value = -9.090909090909
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = -60
$ws.Range("N29").Value = -83.333333333333
$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = -50
$ws.Range("N30").Value = -83.333333333333
